$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 9530089
$ws.Range("J70").Value = 4171.933
$ws.Range("L70").Value = 12515.799
$ws.Range("N70").Value = -13055.799
$ws.Range("H73").Value = 9530089
$ws.Range("J73").Value = 4171.933
$ws.Range("L73").Value = 12515.799
$ws.Range("N73").Value = -14387.799
$ws.Range("H98").Value = 18455
$ws.Range("I98").Value = 21789.691
$ws.Range("J98").Value = 8821.444
$ws.Range("K98").Value = 21789.691
$ws.Range("L98").Value = 8821.444
$ws.Range("M98").Value = -20291.691
$ws.Range("N98").Value = -11817.444
$ws.Range("H122").Value = 18455
$ws.Range("I122").Value = 21789.691
$ws.Range("J122").Value = 8821.444
$ws.Range("K122").Value = 65369.073
$ws.Range("L122").Value = 26464.332
$ws.Range("M122").Value = -62919.073
$ws.Range("N122").Value = -31364.332
$ws.Range("H138").Value = 5916.4
$ws.Range("J138").Value = 6807.25
$ws.Range("L138").Value = 20421.75
$ws.Range("N138").Value = -30701.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 49936.652
$ws.Range("J2").Value = 92826.45
$ws.Range("L2").Value = 92826.45
$ws.Range("N2").Value = -93052.45
$ws.Range("H45").Value = 168031.69
$ws.Range("I45").Value = 240673.78
$ws.Range("J45").Value = 4587
$ws.Range("K45").Value = 240673.78
$ws.Range("L45").Value = 4587
$ws.Range("M45").Value = -240296.78
$ws.Range("N45").Value = -5341
$ws.Range("H74").Value = 1630.3334
$ws.Range("I74").Value = 1526.037
$ws.Range("K74").Value = 1526.037
$ws.Range("M74").Value = -652.037
$ws.Range("H77").Value = 1630.3334
$ws.Range("I77").Value = 1526.037
$ws.Range("K77").Value = 7630.185
$ws.Range("M77").Value = -3262.185
$ws.Range("H116").Value = 49936.652
$ws.Range("J116").Value = 92826.45
$ws.Range("L116").Value = 92826.45
$ws.Range("N116").Value = -97414.45
$ws.Range("H122").Value = 522969.72
$ws.Range("I122").Value = 3778.2632
$ws.Range("K122").Value = 11334.7896
$ws.Range("M122").Value = -8884.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 49936.652
$ws.Range("J3").Value = 92826.45
$ws.Range("L3").Value = 92826.45
$ws.Range("N3").Value = -93054.45
$ws.Range("H20").Value = 2704.5625
$ws.Range("J20").Value = 3838.5
$ws.Range("L20").Value = 3838.5
$ws.Range("N20").Value = -4332.5
$ws.Range("H86").Value = 5676.4194
$ws.Range("I86").Value = 7127.619
$ws.Range("K86").Value = 7127.619
$ws.Range("M86").Value = -6004.619
$ws.Range("H89").Value = 5676.4194
$ws.Range("I89").Value = 7127.619
$ws.Range("K89").Value = 35638.095
$ws.Range("M89").Value = -30022.095
$ws.Range("H96").Value = 19065.818
$ws.Range("I96").Value = 17972.4
$ws.Range("K96").Value = 17972.4
$ws.Range("M96").Value = -15226.4
$ws.Range("H105").Value = 104221.45
$ws.Range("I105").Value = 113943.6
$ws.Range("K105").Value = 113943.6
$ws.Range("M105").Value = -112196.6
$ws.Range("H107").Value = 4692.2
$ws.Range("I107").Value = 6493.6665
$ws.Range("J107").Value = 1990
$ws.Range("K107").Value = 6493.6665
$ws.Range("L107").Value = 1990
$ws.Range("M107").Value = -4573.6665
$ws.Range("N107").Value = -5830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1614.2
$ws.Range("I16").Value = 1681.8
$ws.Range("K16").Value = 1681.8
$ws.Range("M16").Value = -1394.8
$ws.Range("H58").Value = 6835.4595
$ws.Range("I58").Value = 9123.579
$ws.Range("K58").Value = 9123.579
$ws.Range("M58").Value = -8920.579
$ws.Range("H94").Value = 2576.5715
$ws.Range("I94").Value = 4718.4
$ws.Range("J94").Value = 1386.6666
$ws.Range("K94").Value = 4718.4
$ws.Range("L94").Value = 1386.6666
$ws.Range("M94").Value = -4267.4
$ws.Range("N94").Value = -2288.6666
$ws.Range("H105").Value = 85435.67999999999
$ws.Range("I105").Value = 111700.1
$ws.Range("J105").Value = 2265
$ws.Range("K105").Value = 111700.1
$ws.Range("L105").Value = 2265
$ws.Range("M105").Value = -109953.1
$ws.Range("N105").Value = -5759
$ws.Range("H113").Value = 1614.2
$ws.Range("I113").Value = 1681.8
$ws.Range("K113").Value = 1681.8
$ws.Range("M113").Value = 488.2
$ws.Range("H122").Value = 1602.8182
$ws.Range("I122").Value = 1426.25
$ws.Range("J122").Value = 2073.6667
$ws.Range("K122").Value = 4278.75
$ws.Range("L122").Value = 6221.000100000001
$ws.Range("M122").Value = -1828.75
$ws.Range("N122").Value = -11121.0001
$ws.Range("H132").Value = 49064.637
$ws.Range("J132").Value = 68549.836
$ws.Range("L132").Value = 205649.508
$ws.Range("N132").Value = -210709.508
$ws.Range("H136").Value = 6835.4595
$ws.Range("I136").Value = 9123.579
$ws.Range("K136").Value = 27370.737
$ws.Range("M136").Value = -24820.737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 62564.31
$ws.Range("J37").Value = 62564.31
$ws.Range("L37").Value = 187692.93
$ws.Range("N37").Value = -187916.93
$ws.Range("H98").Value = 1091.1
$ws.Range("I98").Value = 1232.4
$ws.Range("K98").Value = 3697.2
$ws.Range("M98").Value = -2199.2
$ws.Range("H113").Value = 1848.6471
$ws.Range("J113").Value = 1883.0625
$ws.Range("L113").Value = 5649.1875
$ws.Range("N113").Value = -9989.1875
$ws.Range("H115").Value = 287.6
$ws.Range("I115").Value = 290.75
$ws.Range("K115").Value = 872.25
$ws.Range("M115").Value = 302.75
$ws.Range("H129").Value = 2170.2354
$ws.Range("I129").Value = 1829.625
$ws.Range("J129").Value = 2473
$ws.Range("K129").Value = 5488.875
$ws.Range("L129").Value = 7419
$ws.Range("M129").Value = -488.875
$ws.Range("N129").Value = -17419
$ws.Range("H132").Value = 16718601
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 16718601
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 150467409
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -150472469
$ws.Range("H134").Value = 28660.908
$ws.Range("I134").Value = 28660.908
$ws.Range("K134").Value = 85982.724
$ws.Range("M134").Value = -80912.724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 147
$ws.Range("I2").Value = 117.583336
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 117.583336
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -4.583336000000003
$ws.Range("N2").Value = -726
$ws.Range("H70").Value = 6436.9443
$ws.Range("I70").Value = 6382.2812
$ws.Range("K70").Value = 6382.2812
$ws.Range("M70").Value = -6112.2812
$ws.Range("H73").Value = 6436.9443
$ws.Range("I73").Value = 6382.2812
$ws.Range("K73").Value = 6382.2812
$ws.Range("M73").Value = -5446.2812
$ws.Range("H102").Value = 22578.572
$ws.Range("I102").Value = 25841.666
$ws.Range("K102").Value = 25841.666
$ws.Range("M102").Value = -24219.666
$ws.Range("H112").Value = 31999.5
$ws.Range("J112").Value = 31999.5
$ws.Range("L112").Value = 31999.5
$ws.Range("N112").Value = -34215.5
$ws.Range("H126").Value = 16216.723
$ws.Range("I126").Value = 20213
$ws.Range("K126").Value = 60639
$ws.Range("M126").Value = -58169

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12526.167
$ws.Range("I22").Value = 17198.273
$ws.Range("J22").Value = 8572.846
$ws.Range("K22").Value = 17198.273
$ws.Range("L22").Value = 8572.846
$ws.Range("M22").Value = -16903.273
$ws.Range("N22").Value = -9162.846
$ws.Range("H27").Value = 12526.167
$ws.Range("I27").Value = 17198.273
$ws.Range("J27").Value = 8572.846
$ws.Range("K27").Value = 17198.273
$ws.Range("L27").Value = 8572.846
$ws.Range("M27").Value = -17091.273
$ws.Range("N27").Value = -8786.846
$ws.Range("H46").Value = 2122.25
$ws.Range("I46").Value = 1194.8
$ws.Range("J46").Value = 2784.7144
$ws.Range("K46").Value = 1194.8
$ws.Range("L46").Value = 2784.7144
$ws.Range("M46").Value = -1006.8
$ws.Range("N46").Value = -3160.7144
$ws.Range("H136").Value = 7331.1924
$ws.Range("I136").Value = 3451.625
$ws.Range("K136").Value = 10354.875
$ws.Range("M136").Value = -7804.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 27001.666
$ws.Range("J21").Value = 32995
$ws.Range("L21").Value = 32995
$ws.Range("N21").Value = -33465
$ws.Range("H32").Value = 20014.5
$ws.Range("J32").Value = 10029
$ws.Range("L32").Value = 10029
$ws.Range("N32").Value = -10663
$ws.Range("H35").Value = 27001.666
$ws.Range("J35").Value = 32995
$ws.Range("L35").Value = 32995
$ws.Range("N35").Value = -33575
$ws.Range("H62").Value = 241629.31
$ws.Range("I62").Value = 444355.16
$ws.Range("J62").Value = 5115.8335
$ws.Range("K62").Value = 444355.16
$ws.Range("L62").Value = 5115.8335
$ws.Range("M62").Value = -443731.16
$ws.Range("N62").Value = -6363.8335
$ws.Range("H65").Value = 241629.31
$ws.Range("I65").Value = 444355.16
$ws.Range("J65").Value = 5115.8335
$ws.Range("K65").Value = 2221775.8
$ws.Range("L65").Value = 25579.1675
$ws.Range("M65").Value = -2218655.8
$ws.Range("N65").Value = -31819.1675
$ws.Range("H136").Value = 3566.5107
$ws.Range("I136").Value = 2881.1282
$ws.Range("K136").Value = 8643.384600000001
$ws.Range("M136").Value = -6093.384600000001
